$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 8724.5
$ws.Range("J10").Value = 8724.5
$ws.Range("L10").Value = 8724.5
$ws.Range("N10").Value = -9310.5
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 200
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = -30
$ws.Range("H13").Value = 2503
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2503
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2503
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -2841
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H80").Value = 750
$ws.Range("J80").Value = 750
$ws.Range("L80").Value = 2250
$ws.Range("N80").Value = -4246
$ws.Range("H83").Value = 750
$ws.Range("J83").Value = 750
$ws.Range("L83").Value = 6750
$ws.Range("N83").Value = -16734
$ws.Range("H107").Value = 1416.1111
$ws.Range("J107").Value = 4000
$ws.Range("L107").Value = 4000
$ws.Range("N107").Value = -7840
$ws.Range("H111").Value = 355
$ws.Range("I111").Value = 353.33334
$ws.Range("J111").Value = 360
$ws.Range("K111").Value = 1060.00002
$ws.Range("L111").Value = 1080
$ws.Range("M111").Value = 2006.99998
$ws.Range("N111").Value = -7214
$ws.Range("H135").Value = 1012.3571
$ws.Range("I135").Value = 859.46155
$ws.Range("K135").Value = 7735.15395
$ws.Range("M135").Value = -5200.15395

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 799.3333
$ws.Range("I12").Value = 448
$ws.Range("J12").Value = 975
$ws.Range("K12").Value = 448
$ws.Range("L12").Value = 975
$ws.Range("M12").Value = -275
$ws.Range("N12").Value = -1321
$ws.Range("H14").Value = 2100.8333
$ws.Range("I14").Value = 1751.25
$ws.Range("J14").Value = 2800
$ws.Range("K14").Value = 1751.25
$ws.Range("L14").Value = 2800
$ws.Range("M14").Value = -1576.25
$ws.Range("N14").Value = -3150
$ws.Range("H45").Value = 2387.5
$ws.Range("J45").Value = 2775
$ws.Range("L45").Value = 2775
$ws.Range("N45").Value = -3529
$ws.Range("H110").Value = 848.375
$ws.Range("I110").Value = 631.3333
$ws.Range("J110").Value = 1499.5
$ws.Range("K110").Value = 631.3333
$ws.Range("L110").Value = 1499.5
$ws.Range("M110").Value = 1413.6667
$ws.Range("N110").Value = -5589.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 26
$ws.Range("I32").Value = 26
$ws.Range("K32").Value = 26
$ws.Range("M32").Value = 358
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 401.33334
$ws.Range("I2").Value = 352
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 352
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -239
$ws.Range("N2").Value = -726
$ws.Range("H3").Value = 2900.6667
$ws.Range("I3").Value = 3351
$ws.Range("K3").Value = 3351
$ws.Range("M3").Value = -3238
$ws.Range("H5").Value = 2280
$ws.Range("I5").Value = 1040
$ws.Range("J5").Value = 2900
$ws.Range("K5").Value = 1040
$ws.Range("L5").Value = 2900
$ws.Range("M5").Value = -928
$ws.Range("N5").Value = -3124
$ws.Range("H11").Value = 2000
$ws.Range("J11").Value = 500
$ws.Range("L11").Value = 500
$ws.Range("N11").Value = -780
$ws.Range("H12").Value = 348.5
$ws.Range("J12").Value = 412
$ws.Range("L12").Value = 412
$ws.Range("N12").Value = -752
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 4801.3335
$ws.Range("I31").Value = 3301.2222
$ws.Range("J31").Value = 7051.5
$ws.Range("K31").Value = 3301.2222
$ws.Range("L31").Value = 7051.5
$ws.Range("M31").Value = -3006.2222
$ws.Range("N31").Value = -7641.5
$ws.Range("H34").Value = 4801.3335
$ws.Range("I34").Value = 3301.2222
$ws.Range("J34").Value = 7051.5
$ws.Range("K34").Value = 3301.2222
$ws.Range("L34").Value = 7051.5
$ws.Range("M34").Value = -3099.2222
$ws.Range("N34").Value = -7455.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H134").Value = 9999.5
$ws.Range("J134").Value = 9999
$ws.Range("L134").Value = 29997
$ws.Range("N134").Value = -35067

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47
$ws.Range("I2").Value = 17.714285
$ws.Range("K2").Value = 106.28571
$ws.Range("M2").Value = 6.714290000000005
$ws.Range("H46").Value = 966
$ws.Range("J46").Value = 966
$ws.Range("L46").Value = 2898
$ws.Range("N46").Value = -3080
$ws.Range("H107").Value = 441.75
$ws.Range("J107").Value = 441.75
$ws.Range("L107").Value = 1325.25
$ws.Range("N107").Value = -5165.25
$ws.Range("H123").Value = 2000
$ws.Range("I123").Value = 2000
$ws.Range("K123").Value = 6000
$ws.Range("M123").Value = -3550
$ws.Range("H129").Value = 1828
$ws.Range("I129").Value = 1424.8
$ws.Range("K129").Value = 4274.4
$ws.Range("M129").Value = 725.6000000000004
$ws.Range("H131").Value = 2086.4285
$ws.Range("I131").Value = 1695
$ws.Range("J131").Value = 2243
$ws.Range("K131").Value = 5085
$ws.Range("L131").Value = 6729
$ws.Range("M131").Value = -45
$ws.Range("N131").Value = -16809
$ws.Range("H132").Value = 1665
$ws.Range("J132").Value = 1995
$ws.Range("L132").Value = 17955
$ws.Range("N132").Value = -23015
$ws.Range("H138").Value = 4690.909
$ws.Range("I138").Value = 4690.909
$ws.Range("K138").Value = 14072.727
$ws.Range("M138").Value = -8932.726999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H21").Value = 3500000
$ws.Range("I21").Value = 3500000
$ws.Range("K21").Value = 3500000
$ws.Range("M21").Value = -3499827
$ws.Range("H30").Value = 3500000
$ws.Range("I30").Value = 3500000
$ws.Range("K30").Value = 3500000
$ws.Range("M30").Value = -3499895
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H122").Value = 17332
$ws.Range("I122").Value = 17332
$ws.Range("K122").Value = 51996
$ws.Range("M122").Value = -49546
$ws.Range("H132").Value = 2588.6667
$ws.Range("I132").Value = 2588.6667
$ws.Range("K132").Value = 7766.000100000001
$ws.Range("M132").Value = -5236.000100000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2333.6
$ws.Range("I16").Value = 2333.6
$ws.Range("K16").Value = 2333.6
$ws.Range("M16").Value = -2163.6
$ws.Range("H132").Value = 4525
$ws.Range("I132").Value = 4525
$ws.Range("K132").Value = 13575
$ws.Range("M132").Value = -11045
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8002.5
$ws.Range("I15").Value = 6006
$ws.Range("K15").Value = 6006
$ws.Range("M15").Value = -5718
$ws.Range("H18").Value = 8002.5
$ws.Range("I18").Value = 6006
$ws.Range("K18").Value = 6006
$ws.Range("M18").Value = -5833
$ws.Range("H20").Value = 22004.5
$ws.Range("I20").Value = 14010
$ws.Range("K20").Value = 14010
$ws.Range("M20").Value = -13770
$ws.Range("H24").Value = 1681669.4
$ws.Range("I24").Value = 2507504.5
$ws.Range("K24").Value = 2507504.5
$ws.Range("M24").Value = -2507274.5
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H31").Value = 9405.666999999999
$ws.Range("H122").Value = 4339.75
$ws.Range("I122").Value = 4786.5
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 14359.5
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -11909.5
$ws.Range("N122").Value = -13898.5
$ws.Range("H128").Value = 53666.668
$ws.Range("J128").Value = 53666.668
$ws.Range("L128").Value = 53666.668
$ws.Range("N128").Value = -63626.668
$ws.Range("H132").Value = 6969.077
$ws.Range("I132").Value = 6240.727
$ws.Range("K132").Value = 18722.181
$ws.Range("M132").Value = -16192.181
